$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Januari 2019")
$ws.Activate()

# --- Window / view bookkeeping (best effort) --------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1

# --- B2: Name field -> literal placeholder text -----------------------------
$ws.Range("B2").Value = "<user>"

# --- Clear the numeric hour entries for these two rows (was 6 / 1) ---------
$ws.Range("B22").ClearContents()
$ws.Range("B28").ClearContents()

# --- Convert colon-style durations to decimal-hour text ---------------------
# Writing "4.5" etc. directly would be auto-parsed as a number; force text
# storage with a leading apostrophe, then restore the original cell format
# (border/fill/number-format) by pasting formats from the neighboring cell
# so the style index is unaffected by the quote-prefix flag.
$ws.Range("B31").Value = "'4.5"
$ws.Range("C31").Copy()
$ws.Range("B31").PasteSpecial(-4122)

$ws.Range("B32").Value = "'2.5"
$ws.Range("C32").Copy()
$ws.Range("B32").PasteSpecial(-4122)

$ws.Range("B33").Value = "'6.5"
$ws.Range("C33").Copy()
$ws.Range("B33").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Replace the hardcoded project total with a live formula ---------------
$ws.Range("B39").Formula = "=SUM(B8:B38)"

# --- Selection as left by the editor -----------------------------------------
$ws.Range("O8").Select()

Write-Host "done"
